$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P2").Value = 0

$ws.Range("H3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0

$ws.Range("H4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0

$ws.Range("P5").Value = 0

$ws.Range("H6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("P6").Value = 0

$ws.Range("H7").Value = 0
$ws.Range("P7").Value = 0

$ws.Range("P8").Value = 0

$ws.Range("P9").Value = 0

$ws.Range("H10").Value = 0
$ws.Range("P10").Value = 0

$ws.Range("H11").Value = 0
$ws.Range("P11").Value = 0

$ws.Range("H12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("P12").Value = 0

$ws.Range("H13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0

$ws.Range("H14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0

$ws.Range("H15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("P15").Value = 0

# Update the view: scroll so column C is the top-left visible column,
# and move the active selection to K15.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("K15").Select()
